# Applies the fixes described in commit "#140, #141, #142 - fixed."
#
#  #140 - "использует лишь" -> "используется лишь" (missing reflexive
#          particle "ся").
#  #141 - "болезненым" -> "болезненным" (missing letter "н").
#  #142 - "их использовать сложнее, чем исходную библиотеку" ->
#          "пользоваться ими сложнее, чем исходной библиотекой"
#          (rewording for correct case agreement).
#
# Plus the two cosmetic/automatic side effects that Word itself performs
# when the author re-types text in the middle of a grammar-checked
# sentence and then saves the document:
#   * the gramStart/gramEnd proof-reading marks around "взаимодействия"
#     disappear because that sentence got re-validated after being
#     edited as a whole;
#   * the hidden "_GoBack" bookmark (Word's "last edit location" marker)
#     moves from the old edit spot (by the "Client" paragraph) to the
#     location of the final edit made in this pass (the "Traces"
#     paragraph).

$d = $word.ActiveDocument

# --- #140: использует -> используется -------------------------------
$d.Content.Find.Execute(
    "использует лишь", $true, $false, $false, $false, $false,
    $true, 1, $false, "используется лишь", 2)

# --- tidy up the "взаимодействия" sentence ---------------------------
# Re-saving the whole (previously gram-checked) sentence as one run
# clears the now-stale gramStart/gramEnd proofErr markers around
# "взаимодействия", matching a fresh grammar pass over the edited text.
$d.Content.Find.Execute(
    "типовых задач взаимодействия с базой данных",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "типовых задач взаимодействия с базой данных", 2)

# --- #141: болезненым -> болезненным ----------------------------------
$d.Content.Find.Execute(
    "болезненым", $true, $false, $false, $false, $false,
    $true, 1, $false, "болезненным", 2)

# --- #142: reword the last sentence about bulky facades ---------------
$d.Content.Find.Execute(
    "их использовать сложнее, чем исходную библиотеку",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "пользоваться ими сложнее, чем исходной библиотекой", 2)

# --- relocate the hidden "_GoBack" bookmark ---------------------------
# Word always keeps a single "_GoBack" bookmark pointing at the location
# of the most recent edit. Since the last textual change made above sits
# right after "...исходной библиотекой" (before the closing period),
# that is where "_GoBack" ends up once the document is saved.
$target = $d.Content
$target.Find.Execute(
    "и пользоваться ими сложнее, чем исходной библиотекой",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
